$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '61.056.78'
$r.Style = 'Normal'
$ws.Range('E2').Value = '  -2.19%  '
$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '2.972.96'
$r.Style = 'Normal'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('E4').Value = '  +0.13%  '
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '592.88'
$r.Style = 'Normal'
$ws.Range('E5').Value = '  +1.55%  '
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '141.86'
$r.Style = 'Normal'
$ws.Range('E6').Value = '  -3.12%  '
$ws.Range('E7').Value = '  +0.01%  '
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '0.513'
$r.Style = 'Normal'
$ws.Range('E8').Value = '  -1.79%  '
$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '2.971.83'
$r.Style = 'Normal'
$ws.Range('E9').Value = '  -1.06%  '
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '0.145'
$r.Style = 'Normal'
$ws.Range('E10').Value = '  -2.22%  '
$ws.Range('E11').Value = '  +4.06%  '
$ws.Range('E12').Value = '  +1.92%  '
$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '0.0000226'
$r.Style = 'Normal'
$ws.Range('E13').Value = '  -1.08%  '
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '33.95'
$r.Style = 'Normal'
$ws.Range('E14').Value = '  -2.02%  '
$ws.Range('E15').Value = '  +1.69%  '
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '3.463.00'
$r.Style = 'Normal'
$ws.Range('E16').Value = '  -0.93%  '
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '61.197.23'
$r.Style = 'Normal'
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('E18').Value = '  -2.82%  '
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '2.969.33'
$r.Style = 'Normal'
$ws.Range('E19').Value = '  -1.14%  '
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '448.25'
$r.Style = 'Normal'
$ws.Range('E20').Value = '  -2.22%  '
$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '14.05'
$r.Style = 'Normal'
$ws.Range('E21').Value = '  +1.24%  '
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '0.677'
$r.Style = 'Normal'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('E23').Value = '  -1.99%  '
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '82.33'
$r.Style = 'Normal'
$ws.Range('E24').Value = '  +2.89%  '
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '2.14'
$r.Style = 'Normal'
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '11.86'
$r.Style = 'Normal'
$ws.Range('E26').Value = '  -2.97%  '
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '10.19'
$r.Style = 'Normal'
$ws.Range('E27').Value = '  +2.16%  '
$ws.Range('E28').Value = '  +0.15%  '
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '0.999'
$r.Style = 'Normal'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  +1.54%  '
$ws.Range('E31').Value = '  -2.24%  '
$ws.Range('E32').Value = '  -3.10%  '
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '27.06'
$r.Style = 'Normal'
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('E35').Value = '  +1.98%  '
$ws.Range('E36').Value = '  -1.25%  '
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '5.73'
$r.Style = 'Normal'
$ws.Range('E37').Value = '  -0.09%  '
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '50.21'
$r.Style = 'Normal'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('E39').Value = '  -3.50%  '
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '8.92'
$r.Style = 'Normal'
$ws.Range('E40').Value = '  +0.30%  '
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '0.120'
$r.Style = 'Normal'
$ws.Range('E41').Value = '  +7.02%  '
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '2.80'
$r.Style = 'Normal'
$ws.Range('E42').Value = '  -4.77%  '
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '387.31'
$r.Style = 'Normal'
$ws.Range('E43').Value = '  -5.47%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '0.0346'
$r.Style = 'Normal'
$ws.Range('E44').Value = '  -1.49%  '
$ws.Range('B45').Value = 'Arweave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '38.43'
$r.Style = 'Normal'
$ws.Range('E45').Value = '  -2.18%  '
$ws.Range('E46').Value = '  -4.59%  '
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '2.684.63'
$r.Style = 'Normal'
$ws.Range('E47').Value = '  -3.19%  '
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '129.48'
$r.Style = 'Normal'
$ws.Range('E48').Value = '  +1.59%  '
$ws.Range('E49').Value = '  +0.14%  '
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '0.107'
$r.Style = 'Normal'
$ws.Range('E50').Value = '  -1.27%  '
$ws.Range('E51').Value = '  -1.13%  '
